$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Pump" test rows (2 and 3) and fill in a new one (row 4) ---
# Row 2: pump run #10
$ws.Range("A2").Value = 10
$ws.Range("C2").Value = "F.0.1.13_1"
$ws.Range("E2").Value = 5.6
$ws.Range("F2").Value = 6.8
$ws.Range("H2").Value = 0.12

# Row 3: pump run #11
$ws.Range("A3").Value = 11
$ws.Range("C3").Value = "F.0.1.13_1"
$ws.Range("E3").Value = 5.6
$ws.Range("F3").Value = 6.8
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 0.12

# Row 4: brand new pump run #12 (previously empty)
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "F.0.1.13_1"
$ws.Range("D4").Value = 1440
$ws.Range("E4").Value = 5.6
$ws.Range("F4").Value = 6.8
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 0.12

# Extend the formatted (bordered/filled) table area down into row 7, matching
# the look of the rows above it.
$ws.Range("A6:H6").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the selected cell shown when the sheet is opened.
$ws.Range("C9").Select()
